$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from existing header cell (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean outlier flag values for rows 2-25, columns F, G, H
$values = @(
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,1),
    @(1,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,0,1),
    @(0,0,0)
)

$row = 2
foreach ($rowValues in $values) {
    $ws.Cells.Item($row, 6).Value = [bool]($rowValues[0])
    $ws.Cells.Item($row, 7).Value = [bool]($rowValues[1])
    $ws.Cells.Item($row, 8).Value = [bool]($rowValues[2])
    $row++
}
